# "Axar Patel" ball-by-ball sheet: the per-innings rows (2-10) in columns
# C:F (runs, balls, fours, sixes) got reshuffled/updated ("updated activity
# till excel form"). Row 1 (headers) and columns A/B (player/team, identical
# on every row) are untouched; row 4 keeps its original values too.
#
# Source cells are text-typed (stored as strings, e.g. t="str"/"21" rather
# than a numeric 21), so each target cell gets NumberFormat "@" (Text)
# applied before the write to keep it a string instead of Excel
# auto-coercing the digits-only literal into a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "C2";  Value = "1"  }
    @{ Cell = "D2";  Value = "4"  }
    @{ Cell = "F2";  Value = "0"  }

    @{ Cell = "C3";  Value = "42" }
    @{ Cell = "D3";  Value = "33" }
    @{ Cell = "E3";  Value = "2"  }
    @{ Cell = "F3";  Value = "3"  }

    @{ Cell = "C5";  Value = "7"  }
    @{ Cell = "E5";  Value = "1"  }

    @{ Cell = "C6";  Value = "21" }
    @{ Cell = "D6";  Value = "5"  }
    @{ Cell = "F6";  Value = "3"  }

    @{ Cell = "C7";  Value = "9"  }
    @{ Cell = "E7";  Value = "1"  }

    @{ Cell = "C8";  Value = "5"  }
    @{ Cell = "D8";  Value = "6"  }
    @{ Cell = "E8";  Value = "0"  }
    @{ Cell = "F8";  Value = "0"  }

    @{ Cell = "C9";  Value = "6"  }
    @{ Cell = "E9";  Value = "0"  }

    @{ Cell = "C10"; Value = "17" }
    @{ Cell = "D10"; Value = "8"  }
    @{ Cell = "F10"; Value = "1"  }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
